# Applies: insert a new "2022-Q3" sheet (with fund holdings detail) right
# after the "总计" (summary) sheet, and add a corresponding summary row in
# the "总计" sheet. All pre-existing quarter sheets shift right by one
# position but keep their own data unchanged.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, [string]$text) {
    # Force a string value onto a cell even when the text looks numeric
    # (keeps leading zeros / trailing zeros such as "004738" or "3.70"),
    # then strip the quote-prefix style artifact the apostrophe leaves
    # behind so the cell's style matches a normal, unformatted cell.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Summary sheet "总计": insert a new row for 2022-Q3 at the top of the
#    data (row 2), pushing all existing quarters down by one row, then
#    renumber the sequential index column (A). Rows are shifted by
#    copying values bottom-up (rather than Rows.Insert, which drags the
#    header row's bold style into the new blank row) so formatting of
#    the untouched rows is preserved exactly.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Cells.Item($dst, 1).Value = $summary.Cells.Item($r, 1).Value()
    $summary.Cells.Item($dst, 2).Value = $summary.Cells.Item($r, 2).Value()
    $summary.Cells.Item($dst, 3).Value = $summary.Cells.Item($r, 3).Value()
    $summary.Cells.Item($dst, 4).Value = $summary.Cells.Item($r, 4).Value()
}

# Row 9 is brand new territory (sheet used to stop at row 8) - copy A8's
# style onto A9 so the index column stays consistently formatted.
$summary.Cells.Item(8, 1).Copy()
$summary.Cells.Item(9, 1).PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
Set-TextCell $summary.Cells.Item(2, 2) "2022-Q3"
$summary.Cells.Item(2, 3).Value = 11
$summary.Cells.Item(2, 4).Value = 1.69

for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" detail sheet right after "总计" by copying
#    the existing "2022-Q2" sheet (so header/row styling matches the
#    other quarter sheets exactly) and overwriting its data with the
#    2022-Q3 fund holdings.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Copy([System.Reflection.Missing]::Value, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3Data = @(
    @("159611", "广发中证全指电力ETF",           "16.62", "99.58", "4.99", "0.8293", 7),
    @("004958", "圆信永丰优享生活灵活配置混合",     "21.77", "80.56", "1.68", "0.3657", 10),
    @("159625", "嘉实国证绿色电力ETF",           "3.41",  "98.77", "3.70", "0.1262", 5),
    @("004823", "上投摩根安裕回报混合A",          "5.35",  "25.71", "1.93", "0.1033", 3),
    @("004824", "上投摩根安裕回报混合C",          "4.91",  "25.71", "1.93", "0.0948", 3),
    @("561700", "博时中证全指电力公用事业ETF",     "1.31",  "98.79", "4.95", "0.0648", 7),
    @("561560", "华泰柏瑞中证全指电力公用事业ETF", "0.90",  "98.29", "4.98", "0.0448", 7),
    @("519615", "银河君尚灵活配置混合I",          "3.59",  "35.36", "0.73", "0.0262", 6),
    @("562350", "银华中证全指电力公用事业ETF",     "0.46",  "97.99", "4.91", "0.0226", 7),
    @("519613", "银河君尚灵活配置混合A",          "2.10",  "35.36", "0.73", "0.0153", 6),
    @("519614", "银河君尚灵活配置混合C",          "0.17",  "35.36", "0.73", "0.0012", 6)
)

$row = 2
foreach ($fund in $q3Data) {
    Set-TextCell $q3.Cells.Item($row, 2) $fund[0]
    Set-TextCell $q3.Cells.Item($row, 3) $fund[1]
    Set-TextCell $q3.Cells.Item($row, 4) $fund[2]
    Set-TextCell $q3.Cells.Item($row, 5) $fund[3]
    Set-TextCell $q3.Cells.Item($row, 6) $fund[4]
    Set-TextCell $q3.Cells.Item($row, 7) $fund[5]
    $q3.Cells.Item($row, 8).Value = $fund[6]
    $row++
}
